$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.007715731536047051
$ws.Range("C2").Value = 0.01072926463919976
$ws.Range("D2").Value = 0.01313059040095384
$ws.Range("E2").Value = 0.01533374246478889

$ws.Range("B3").Value = 2.960985194873507
$ws.Range("C3").Value = 5.278654477475951
$ws.Range("D3").Value = 7.920378794277832
$ws.Range("E3").Value = 6.154797652350622

$ws.Range("B4").Value = -0.01923241907602552
$ws.Range("C4").Value = -0.02211886608468998
$ws.Range("D4").Value = -0.02496364126227692
$ws.Range("E4").ClearContents()

$ws.Range("B5").Value = -7.248605109103126
$ws.Range("C5").Value = -8.497494301771022
$ws.Range("D5").Value = -8.279145604947555
$ws.Range("E5").ClearContents()

$ws.Range("B6").Value = 0.008306192841956505
$ws.Range("C6").Value = 0.004700624964746629
$ws.Range("D6").Value = 0.0008618560449260675
$ws.Range("E6").ClearContents()

$ws.Range("B7").Value = 2.19958781331111
$ws.Range("C7").Value = 1.449288110650679
$ws.Range("D7").Value = 0.2842886919970967
$ws.Range("E7").ClearContents()

$ws.Range("B8").Value = 0.008357527046492438
$ws.Range("C8").Value = 0.01192195670660364
$ws.Range("D8").Value = 0.01434129409155893
$ws.Range("E8").Value = 0.01681369197704382

$ws.Range("B9").Value = 3.423039717211485
$ws.Range("C9").Value = 5.899705947516384
$ws.Range("D9").Value = 6.631425257262459
$ws.Range("E9").Value = 6.430676064501979

$ws.Range("B10").Value = -0.01895423313014715
$ws.Range("C10").Value = -0.02017157522180509
$ws.Range("D10").Value = -0.02473949831779778
$ws.Range("E10").ClearContents()

$ws.Range("B11").Value = -6.9755409058152
$ws.Range("C11").Value = -7.944255618650784
$ws.Range("D11").Value = -7.219447958733066
$ws.Range("E11").ClearContents()

$ws.Range("B12").Value = 0.007791251779627479
$ws.Range("C12").Value = 0.00225479725540151
$ws.Range("D12").ClearContents()
$ws.Range("E12").ClearContents()

$ws.Range("B13").Value = 2.263420962238378
$ws.Range("C13").Value = 0.7330567005423656
$ws.Range("D13").ClearContents()
$ws.Range("E13").ClearContents()

$ws.Range("B14").Value = 0.009041690451927127
$ws.Range("C14").Value = 0.01263615281712632
$ws.Range("D14").Value = 0.01594696925871424
$ws.Range("E14").Value = 0.01836798765851275

$ws.Range("B15").Value = 3.479376881707805
$ws.Range("C15").Value = 5.834334341290351
$ws.Range("D15").Value = 6.492683760118034
$ws.Range("E15").Value = 7.104735858558067

$ws.Range("B16").Value = -0.01489975316747699
$ws.Range("C16").Value = -0.0174307244100564
$ws.Range("D16").Value = -0.02030460099731261
$ws.Range("E16").ClearContents()

$ws.Range("B17").Value = -5.535363787294505
$ws.Range("C17").Value = -6.53827035388502
$ws.Range("D17").Value = -6.086614299188402
$ws.Range("E17").ClearContents()

$ws.Range("B18").Value = 0.003228875245632952
$ws.Range("C18").Value = -0.001599139166869575
$ws.Range("D18").ClearContents()
$ws.Range("E18").ClearContents()

$ws.Range("B19").Value = 0.9402343593677677
$ws.Range("C19").Value = -0.4787012132139899
$ws.Range("D19").ClearContents()
$ws.Range("E19").ClearContents()

$ws.Range("B20").Value = 0.01206491642932434
$ws.Range("C20").Value = 0.01610843657722701
$ws.Range("D20").Value = 0.01945519334682248
$ws.Range("E20").Value = 0.02070905493725941

$ws.Range("B21").Value = 4.534244893993776
$ws.Range("C21").Value = 7.321129523134057
$ws.Range("D21").Value = 7.606929816504106
$ws.Range("E21").Value = 7.922102879841858

$ws.Range("B22").Value = -0.01426144422137027
$ws.Range("C22").Value = -0.01575796625718798
$ws.Range("D22").Value = -0.01791234271276543
$ws.Range("E22").ClearContents()

$ws.Range("B23").Value = -5.027115627469286
$ws.Range("C23").Value = -5.801771728919297
$ws.Range("D23").Value = -5.505493227295577
$ws.Range("E23").ClearContents()

$ws.Range("B24").Value = -0.0007985145805884673
$ws.Range("C24").Value = -0.007556686967392715
$ws.Range("D24").ClearContents()
$ws.Range("E24").ClearContents()

$ws.Range("B25").Value = -0.2135405384522041
$ws.Range("C25").Value = -2.094864827514199
$ws.Range("D25").ClearContents()
$ws.Range("E25").ClearContents()
